$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. "Tomasz Kisiel" paragraph: drop the pl-PL language run/paragraph formatting ---
$pName = $d.Paragraphs(1)
$pName.Range.InsertXML('<w:p ' + $wNs + '><w:r><w:t>Tomasz Kisiel</w:t></w:r></w:p>') | Out-Null

# --- 2. "27-10-2024" paragraph: drop the pl-PL language run/paragraph formatting ---
$pDate = $d.Paragraphs(2)
$pDate.Range.InsertXML('<w:p ' + $wNs + '><w:r><w:t>27-10-2024</w:t></w:r></w:p>') | Out-Null

# --- 3. following empty paragraph: drop the pl-PL paragraph formatting entirely ---
$pBlank = $d.Paragraphs(3)
$pBlank.Range.InsertXML('<w:p ' + $wNs + '/>') | Out-Null

# --- 4. Task 2 paragraph: merge the two runs describing average-weight sorting into one run ---
$mergedText = "Calculate average weight for each type of fish. Sort them in decreasing order. Present only 50 the heaviest fish."
$d.Content.Find.Execute($mergedText, $true, $false, $false, $false, $false, $true, 1, $false, $mergedText, 2) | Out-Null

# --- 5. Task 5 "Code" paragraph: change the example code from 3SME to 5BRE and split across 3 runs ---
$quoteOpen = [char]8220
$quoteClose = [char]8221
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Combine the column from task 4*") {
        $xml = '<w:p ' + $wNs + '>' +
               '<w:r><w:t xml:space="preserve">Combine the column from task 4 and task 5 in order to create a new column Code in a form of ' + $quoteOpen + '3SME' + $quoteClose + '. Count how </w:t></w:r>' +
               '<w:r><w:t>many observations will have ' + $quoteOpen + '5BR</w:t></w:r>' +
               '<w:r><w:t>E' + $quoteClose + ' code.</w:t></w:r>' +
               '</w:p>'
        $para.Range.InsertXML($xml) | Out-Null
        break
    }
}

# --- 6. Task 8: merge the split "Save the complete..." runs, and move the _GoBack bookmark
#         out into its own new paragraph (consuming two of the trailing blank paragraphs) ---
$saveIdx = -1
$idx = 0
foreach ($para in $d.Paragraphs) {
    $idx = $idx + 1
    if ($para.Range.Text -like "Save the complete fish table wit*") {
        $saveIdx = $idx
        break
    }
}
$pSave = $d.Paragraphs($saveIdx)
$pBlank2 = $d.Paragraphs($saveIdx + 2)
$startPos = $pSave.Range.Start
$endPos = $pBlank2.Range.End
$target = $d.Range($startPos, $endPos)
$xml = '<w:p ' + $wNs + '><w:r><w:t>Save the complete fish table with 159 observations and new columns into an excel file.</w:t></w:r></w:p>' +
       '<w:p ' + $wNs + '/>' +
       '<w:p ' + $wNs + '><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$target.InsertXML($xml) | Out-Null
